# Auto-generated script applying the cryptos.xlsx update
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '29.918.34'
$ws.Range('E2').Value = '  +0.06%  '
$ws.Range('D3').Value = '1.876.37'
$ws.Range('E3').Value = '  -0.61%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '0.7424'
$ws.Range('E5').Value = '  -3.87%  '
$ws.Range('D6').Value = '242.57'
$ws.Range('E6').Value = '  -0.19%  '
$ws.Range('D8').Value = '0.3154'
$ws.Range('E8').Value = '  +1.50%  '
$ws.Range('D9').Value = '0.07261'
$ws.Range('E9').Value = '  +1.09%  '
$ws.Range('E10').Value = '  -3.54%  '
$ws.Range('D11').Value = '0.08401'
$ws.Range('E11').Value = '  -2.00%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.909.42'
$ws.Range('E12').Value = '  -0.48%  '
$ws.Range('B13').Value = 'Polygon'
$ws.Range('C13').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D13').Value = '0.7522'
$ws.Range('E13').Value = '  -1.50%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').Value = '5.438'
$ws.Range('E14').Value = '  +1.19%  '
$ws.Range('D15').Value = '92.55'
$ws.Range('E15').Value = '  -1.27%  '
$ws.Range('D16').Value = '29.943.75'
$ws.Range('E16').Value = '  +0.09%  '
$ws.Range('D17').Value = '6.089'
$ws.Range('E17').Value = '  -1.71%  '
$ws.Range('D18').Value = '247.39'
$ws.Range('E18').Value = '  +1.25%  '
$ws.Range('E19').Value = '  -1.38%  '
$ws.Range('D20').Value = '0.000007870'
$ws.Range('E20').Value = '  +0.74%  '
$ws.Range('E21').Value = '  +0.20%  '
$ws.Range('D22').Value = '2.127.40'
$ws.Range('E22').Value = '  -1.90%  '
$ws.Range('D23').Value = '8.051'
$ws.Range('E23').Value = '  +1.01%  '
$ws.Range('E24').Value = '  -0.06%  '
$ws.Range('D25').Value = '0.1561'
$ws.Range('E25').Value = '  -5.21%  '
$ws.Range('D26').Value = '9.276'
$ws.Range('E26').Value = '  -1.04%  '
$ws.Range('D27').Value = '165.22'
$ws.Range('E27').Value = '  +2.15%  '
$ws.Range('D28').Value = '18.64'
$ws.Range('E28').Value = '  -0.50%  '
$ws.Range('E29').Value = '  +0.27%  '
$ws.Range('D30').Value = '1.514'
$ws.Range('E30').Value = '  +4.79%  '
$ws.Range('D31').Value = '4.606'
$ws.Range('E31').Value = '  +1.90%  '
$ws.Range('D32').Value = '1.539'
$ws.Range('E32').Value = '  +0.16%  '
$ws.Range('D33').Value = '4.293'
$ws.Range('E33').Value = '  +4.60%  '
$ws.Range('D34').Value = '0.05346'
$ws.Range('E34').Value = '  -1.74%  '
$ws.Range('E35').Value = '  -0.17%  '
$ws.Range('D36').Value = '0.7535'
$ws.Range('E36').Value = '  +1.07%  '
$ws.Range('E37').Value = '  -0.49%  '
$ws.Range('D38').Value = '2.691'
$ws.Range('E38').Value = '  -0.03%  '
$ws.Range('D39').Value = '0.01967'
$ws.Range('E39').Value = '  +0.02%  '
$ws.Range('D40').Value = '2.762'
$ws.Range('E40').Value = '  -0.67%  '
$ws.Range('D41').Value = '0.4542'
$ws.Range('E41').Value = '  +1.78%  '
$ws.Range('D42').Value = '1.112.82'
$ws.Range('E42').Value = '  +0.39%  '
$ws.Range('D43').Value = '6.052'
$ws.Range('E43').Value = '  -0.55%  '
$ws.Range('D44').Value = '72.66'
$ws.Range('E44').Value = '  -0.79%  '
$ws.Range('E45').Value = '  +0.44%  '
$ws.Range('E46').Value = '  +0.12%  '
$ws.Range('D47').Value = '103.43'
$ws.Range('E47').Value = '  +0.59%  '
$ws.Range('B48').Value = 'Aptos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D48').Value = '7.628'
$ws.Range('E48').Value = '  +0.03%  '
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').Value = '1.858'
$ws.Range('E49').Value = '  -0.67%  '
$ws.Range('D50').Value = '2.025.02'
$ws.Range('E50').Value = '  -3.22%  '
$ws.Range('D51').Value = '2.903'
$ws.Range('E51').Value = '  -2.79%  '
